$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(45958, 8630.47329661029, 8024.00919333907, 12068, 3912.99505538061, -5.45815630334649),
    @(45959, 8042.2835020355, 7574.51201279366, 3620, 3981.91642701909, 330.684518325531),
    @(45960, 7622.40156379609, 7015.66280248329, 3620, 3710.87082210575, 296.10556769121),
    @(45961, 7048.05934848961, 5804.86632320966, 3620, 3317.29761372859, 229.25683070576),
    @(45962, 2061.768603528, 3379.9141446472, 11404, 3230.44520232422, -199.735027209524),
    @(45963, 2557.16655997904, 3913.55961206506, 11404, 3797.24632547564, -153.883085935804),
    @(45964, 8571.71877057272, 7182.3872148972, 11404, 4686.02043775852, 19.3503188606548),
    @(45965, 8571.71877057272, 7518.21337313226, 11404, 4686.02043775852, 33.3430754537825),
    @(45966, 8571.71877057272, 7388.28237988266, 11404, 4686.02043775852, 27.9292840683825),
    @(45967, 8580.99785695336, 7443.82207224924, 11404, 4688.10953166867, 30.3304834965796),
    @(45968, 8585.92244945998, 6860.20489713919, 11404, 4719.63295453454, 7.32657715307187),
    @(45969, 3147.97579797074, 4906.35167786925, 11404, 4404.40809265882, -87.2183428946638),
    @(45970, 3012.17181266658, 5042.55828274956, 11404, 4396.26845485097, -81.8822192666447),
    @(45971, 9340.71276866398, 8370.94141257818, 11404, 5362.49965599072, 97.0600445237042)
)

$startRow = 2
for ($i = 0; $i -lt $data.Length; $i++) {
    $r = $startRow + $i
    $row = $data[$i]
    for ($c = 0; $c -lt $row.Length; $c++) {
        $ws.Cells.Item($r, $c + 1).Value2 = $row[$c]
    }
}
